# Updated league table for GW21.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add this gameweek's scores for row 18 (Week 20)
$ws.Range("B18").Value = 104
$ws.Range("C18").Value = 127
$ws.Range("D18").Value = 132

# Update the selection to match the author's saved cursor position
$ws.Range("D19").Select()
